# Add team record columns (Wins, Losses, Ties) to the OAK_2000 sheet.
# Mirrors the commit: "Added team record to data" - appends AD/AE/AF
# with constant team totals (91 wins, 70 losses, 0 ties) for every player row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 46

# --- Header row (row 1): copy the existing header style (from AC1) so the
#     new headers match the bold/centered/bordered look of the other headers.
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("AD1").Value = "Wins"

$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AE1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("AE1").Value = "Losses"

$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AF1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("AF1").Value = "Ties"

$excel.CutCopyMode = 0

# --- Data rows (2..46): constant team record values for every player.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 91   # AD -> Wins
    $ws.Cells.Item($r, 31).Value = 70   # AE -> Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF -> Ties
}

Write-Host "Added Wins/Losses/Ties columns (AD:AF) for rows 1-$lastRow"
